$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Sun Oct 13 00:07:34 EDT 2024"
$ws.Range("B3").Value = "Sun Oct 13 00:07:48 EDT 2024"
$ws.Range("B4").Value = "Sun Oct 13 00:08:01 EDT 2024"
